# Fix: added clg, dep in subjects bulk
# Adds two new header columns ("Department ID" in D1, "College ID" in E1)
# to the subjects bulk-upload template, matching the bold header style
# already used by A1:C1.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# D1 previously existed as a styled-but-empty cell; give it its header text.
$ws.Range("D1").Value = "Department ID"
$ws.Range("E1").Value = "College ID"

# Match the bold header formatting used by the rest of row 1.
$ws.Range("D1:E1").Font.Bold = $true
